$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.659.53"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").Value = "1.599.79"
$ws.Range("E3").Value = "  +0.57%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.38"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.07%  "
$ws.Range("E6").Value = "  +1.10%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.0618"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("E10").Value = "  -1.03%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0837"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "1.823.69"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("D13").Value = "1.600.92"
$ws.Range("E13").Value = "  +0.64%  "
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("E15").Value = "  -0.79%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.82"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.10%  "
$ws.Range("D17").Value = "26.645.18"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "208.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.97"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.27%  "
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("E24").Value = "  -0.61%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.58"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.81%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.657"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").Value = "1.283.42"
$ws.Range("E35").Value = "  -2.65%  "
$ws.Range("E36").Value = "  +1.57%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.49"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("E39").Value = "  +1.84%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +1.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.21"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.55%  "
$ws.Range("E43").Value = "  -0.46%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.97"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.919"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.50%  "
$ws.Range("D46").Value = "1.736.27"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "89.78"
$ws.Range("D47").ClearFormats()
$ws.Range("E48").Value = "  -0.70%  "
$ws.Range("E49").Value = "  +3.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0505"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.14%  "
